$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-19 (only changed cells, per diff) ---
$ws.Range("B2").Value = "NSE:EUROTEXIND"
$ws.Range("C2").Value = "NSE:63MOONS"
$ws.Range("E2").Value = "NSE:AARTIIND"
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = "NSE:ABFRL"
$ws.Range("E3").Value = "NSE:CANBK"
$ws.Range("B4").ClearContents()
$ws.Range("C4").Value = "NSE:APEX"
$ws.Range("E4").Value = "NSE:FEDERALBNK"
$ws.Range("B5").ClearContents()
$ws.Range("C5").Value = "NSE:APOLLOHOSP"
$ws.Range("E5").Value = "NSE:INDUSINDBK"
$ws.Range("B6").ClearContents()
$ws.Range("C6").Value = "NSE:ASHOKAMET"
$ws.Range("E6").Value = "NSE:INDUSTOWER"
$ws.Range("B7").ClearContents()
$ws.Range("C7").Value = "NSE:AUTOBEES"
$ws.Range("E7").Value = "NSE:LICHSGFIN"
$ws.Range("B8").ClearContents()
$ws.Range("C8").Value = "NSE:AXSENSEX"
$ws.Range("E8").Value = "NSE:NTPC"
$ws.Range("B9").ClearContents()
$ws.Range("C9").Value = "NSE:BAJAJ-AUTO"
$ws.Range("E9").Value = "NSE:PEL"
$ws.Range("B10").ClearContents()
$ws.Range("C10").Value = "NSE:BATAINDIA"
$ws.Range("E10").Value = "NSE:POWERGRID"
$ws.Range("B11").ClearContents()
$ws.Range("C11").Value = "NSE:BLS"
$ws.Range("B12").ClearContents()
$ws.Range("C12").Value = "NSE:BPCL"
$ws.Range("B13").ClearContents()
$ws.Range("C13").Value = "NSE:CONSUMBEES"
$ws.Range("B14").ClearContents()
$ws.Range("C14").Value = "NSE:DABUR"
$ws.Range("C15").Value = "NSE:DALBHARAT"
$ws.Range("C16").Value = "NSE:DANGEE"
$ws.Range("C17").Value = "NSE:ESG"
$ws.Range("C18").Value = "NSE:GODREJCP"
$ws.Range("C19").Value = "NSE:GODREJIND"

# --- Append new rows 20-43, cloning row-19 style for column A, then set values ---
$ws.Range("A19").Copy($ws.Range("A20"))
$ws.Range("A20").Value = 18
$ws.Range("C20").Value = "NSE:GRASIM"
$ws.Range("A19").Copy($ws.Range("A21"))
$ws.Range("A21").Value = 19
$ws.Range("C21").Value = "NSE:HCL-INSYS"
$ws.Range("A19").Copy($ws.Range("A22"))
$ws.Range("A22").Value = 20
$ws.Range("C22").Value = "NSE:HEROMOTOCO"
$ws.Range("A19").Copy($ws.Range("A23"))
$ws.Range("A23").Value = 21
$ws.Range("C23").Value = "NSE:HINDPETRO"
$ws.Range("A19").Copy($ws.Range("A24"))
$ws.Range("A24").Value = 22
$ws.Range("C24").Value = "NSE:HTMEDIA"
$ws.Range("A19").Copy($ws.Range("A25"))
$ws.Range("A25").Value = 23
$ws.Range("C25").Value = "NSE:IOC"
$ws.Range("A19").Copy($ws.Range("A26"))
$ws.Range("A26").Value = 24
$ws.Range("C26").Value = "NSE:JAYSREETEA"
$ws.Range("A19").Copy($ws.Range("A27"))
$ws.Range("A27").Value = 25
$ws.Range("C27").Value = "NSE:JTLIND"
$ws.Range("A19").Copy($ws.Range("A28"))
$ws.Range("A28").Value = 26
$ws.Range("C28").Value = "NSE:KAJARIACER"
$ws.Range("A19").Copy($ws.Range("A29"))
$ws.Range("A29").Value = 27
$ws.Range("C29").Value = "NSE:KIRIINDUS"
$ws.Range("A19").Copy($ws.Range("A30"))
$ws.Range("A30").Value = 28
$ws.Range("C30").Value = "NSE:KIRLPNU"
$ws.Range("A19").Copy($ws.Range("A31"))
$ws.Range("A31").Value = 29
$ws.Range("C31").Value = "NSE:LICNETFN50"
$ws.Range("A19").Copy($ws.Range("A32"))
$ws.Range("A32").Value = 30
$ws.Range("C32").Value = "NSE:LT"
$ws.Range("A19").Copy($ws.Range("A33"))
$ws.Range("A33").Value = 31
$ws.Range("C33").Value = "NSE:LUMAXIND"
$ws.Range("A19").Copy($ws.Range("A34"))
$ws.Range("A34").Value = 32
$ws.Range("C34").Value = "NSE:MAANALU"
$ws.Range("A19").Copy($ws.Range("A35"))
$ws.Range("A35").Value = 33
$ws.Range("C35").Value = "NSE:MARUTI"
$ws.Range("A19").Copy($ws.Range("A36"))
$ws.Range("A36").Value = 34
$ws.Range("C36").Value = "NSE:MONTECARLO"
$ws.Range("A19").Copy($ws.Range("A37"))
$ws.Range("A37").Value = 35
$ws.Range("C37").Value = "NSE:NAVINFLUOR"
$ws.Range("A19").Copy($ws.Range("A38"))
$ws.Range("A38").Value = 36
$ws.Range("C38").Value = "NSE:NPBET"
$ws.Range("A19").Copy($ws.Range("A39"))
$ws.Range("A39").Value = 37
$ws.Range("C39").Value = "NSE:NSIL"
$ws.Range("A19").Copy($ws.Range("A40"))
$ws.Range("A40").Value = 38
$ws.Range("C40").Value = "NSE:PALASHSECU"
$ws.Range("A19").Copy($ws.Range("A41"))
$ws.Range("A41").Value = 39
$ws.Range("C41").Value = "NSE:PCJEWELLER"
$ws.Range("A19").Copy($ws.Range("A42"))
$ws.Range("A42").Value = 40
$ws.Range("C42").Value = "NSE:PIDILITIND"
$ws.Range("A19").Copy($ws.Range("A43"))
$ws.Range("A43").Value = 41
$ws.Range("C43").Value = "NSE:RELIANCE"
